$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.029.25'
$ws.Range('D3').Value = '1.653.07'
$ws.Range('E3').Value = '  +3.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('E6').Value = '  +1.94%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +1.75%  '
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('E10').Value = '  +4.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0869'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').Value = '1.886.23'
$ws.Range('E12').Value = '  +3.51%  '
$ws.Range('D13').Value = '1.657.16'
$ws.Range('E13').Value = '  +3.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.20%  '
$ws.Range('E16').Value = '  +2.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '239.39'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.38%  '
$ws.Range('D18').Value = '27.034.85'
$ws.Range('E18').Value = '  +2.66%  '
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  +4.19%  '
$ws.Range('E23').Value = '  +3.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.96%  '
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.83'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0499'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.94%  '
$ws.Range('E31').Value = '  +1.57%  '
$ws.Range('E32').Value = '  +3.08%  '
$ws.Range('D33').Value = '1.513.57'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('E34').Value = '  +5.16%  '
$ws.Range('E35').Value = '  +8.90%  '
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.580'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.891'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.03%  '
$ws.Range('E39').Value = '  +3.23%  '
$ws.Range('E40').Value = '  +3.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.26'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = '1.793.74'
$ws.Range('E44').Value = '  +3.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.776'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.918'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.79'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.56%  '
$ws.Range('E48').Value = '  +2.67%  '
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0977'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.08%  '
